$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "HAC", 2025, 1, 1.5267),
    @(3, "HAC", 2025, 2, 1.4794),
    @(4, "HAC", 2025, 3, 1.5908000000000002),
    @(5, "HAC", 2025, 4, 1.3506),
    @(6, "HAC", 2025, 5, 1.5129),
    @(7, "HAC", 2025, 6, 1.3161),
    @(8, "HAC", 2025, 7, 1.6600999999999997),
    @(9, "HAC", 2025, 8, 1.4393),
    @(10, "HAC", 2025, 9, 1.3805),
    @(11, "HAC", 2025, 10, 1.7260000000000002),
    @(12, "HEM", 2025, 1, 1.4035),
    @(13, "HEM", 2025, 2, 1.3419),
    @(14, "HEM", 2025, 3, 1.5225),
    @(15, "HEM", 2025, 4, 1.3316),
    @(16, "HEM", 2025, 5, 1.2706),
    @(17, "HEM", 2025, 6, 1.4464),
    @(18, "HEM", 2025, 7, 1.3363),
    @(19, "HEM", 2025, 8, 1.2408),
    @(20, "HEM", 2025, 9, 1.5421),
    @(21, "HEM", 2025, 10, 1.2747),
    @(22, "HIJPII", 2025, 1, 1.4443),
    @(23, "HIJPII", 2025, 2, 1.3426),
    @(24, "HIJPII", 2025, 3, 1.2255000000000003),
    @(25, "HIJPII", 2025, 4, 1.158),
    @(26, "HIJPII", 2025, 5, 1.1414),
    @(27, "HIJPII", 2025, 6, 0.9257),
    @(28, "HIJPII", 2025, 7, 1.2025),
    @(29, "HIJPII", 2025, 8, 1.3008),
    @(30, "HIJPII", 2025, 9, 1.5947),
    @(31, "HIJPII", 2025, 10, 1.1206),
    @(32, "HJK", 2025, 1, 1.5392),
    @(33, "HJK", 2025, 2, 1.7159),
    @(34, "HJK", 2025, 3, 1.5151),
    @(35, "HJK", 2025, 4, 1.4668),
    @(36, "HJK", 2025, 5, 1.3524999999999998),
    @(37, "HJK", 2025, 6, 1.0959),
    @(38, "HJK", 2025, 7, 1.449),
    @(39, "HJK", 2025, 8, 1.3163),
    @(40, "HJK", 2025, 9, 1.2914),
    @(41, "HJK", 2025, 10, 1.1885),
    @(42, "HJXXIII", 2025, 1, 1.4294),
    @(43, "HJXXIII", 2025, 2, 1.4434),
    @(44, "HJXXIII", 2025, 3, 1.6155),
    @(45, "HJXXIII", 2025, 4, 1.3379),
    @(46, "HJXXIII", 2025, 5, 1.4978),
    @(47, "HJXXIII", 2025, 6, 1.4668),
    @(48, "HJXXIII", 2025, 7, 1.5574),
    @(49, "HJXXIII", 2025, 8, 1.5182),
    @(50, "HJXXIII", 2025, 9, 1.5721),
    @(51, "HJXXIII", 2025, 10, 1.4951),
    @(52, "HRAD", 2025, 1, 1.2246),
    @(53, "HRAD", 2025, 2, 1.0594),
    @(54, "HRAD", 2025, 3, 1.2125),
    @(55, "HRAD", 2025, 4, 0.998),
    @(56, "HRAD", 2025, 5, 1.148),
    @(57, "HRAD", 2025, 6, 1.1968),
    @(58, "HRAD", 2025, 7, 0.9564),
    @(59, "HRAD", 2025, 8, 1.0831),
    @(60, "HRAD", 2025, 9, 0.9619),
    @(61, "HRAD", 2025, 10, 0.8774),
    @(62, "HRBJA", 2025, 1, 1.4969),
    @(63, "HRBJA", 2025, 2, 1.3492),
    @(64, "HRBJA", 2025, 3, 1.4684),
    @(65, "HRBJA", 2025, 4, 1.3237),
    @(66, "HRBJA", 2025, 5, 1.3811),
    @(67, "HRBJA", 2025, 6, 1.4875),
    @(68, "HRBJA", 2025, 7, 1.4807),
    @(69, "HRBJA", 2025, 8, 1.4999),
    @(70, "HRBJA", 2025, 9, 1.783),
    @(71, "HRBJA", 2025, 10, 1.5088),
    @(72, "HRJP", 2025, 1, 1.797),
    @(73, "HRJP", 2025, 2, 1.6135),
    @(74, "HRJP", 2025, 3, 1.6475),
    @(75, "HRJP", 2025, 4, 1.7361),
    @(76, "HRJP", 2025, 5, 1.5634),
    @(77, "HRJP", 2025, 6, 1.5819),
    @(78, "HRJP", 2025, 7, 1.579),
    @(79, "HRJP", 2025, 8, 1.4975),
    @(80, "HRJP", 2025, 9, 1.6762),
    @(81, "HRJP", 2025, 10, 1.4729),
    @(82, "MOV", 2025, 1, 1.3053),
    @(83, "MOV", 2025, 2, 1.2967),
    @(84, "MOV", 2025, 3, 1.4614000000000003),
    @(85, "MOV", 2025, 4, 1.2698),
    @(86, "MOV", 2025, 5, 1.3283999999999998),
    @(87, "MOV", 2025, 6, 1.4484),
    @(88, "MOV", 2025, 7, 1.4048),
    @(89, "MOV", 2025, 8, 1.3474),
    @(90, "MOV", 2025, 9, 1.2602),
    @(91, "MOV", 2025, 10, 1.333)
)

# Rows 83-91 are brand new; copy formatting from the last existing data row (row 82) down to them first.
$ws.Range("A82:D82").Copy()
$ws.Range("A83:D91").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($row in $data) {
    $r = $row[0]
    $unit = $row[1]
    $year = $row[2]
    $month = $row[3]
    $val = $row[4]

    $ws.Cells.Item($r, 1).Value = $unit
    $ws.Cells.Item($r, 2).Value = $year
    $ws.Cells.Item($r, 3).Value = $month
    $ws.Cells.Item($r, 4).Value = $val
}

$ws.Range("A92:XFD93").Select()
